# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    # Leading apostrophe forces Excel to store the value as literal text
    # even when it looks numeric (e.g. "9.27"); reset the style afterwards
    # so no stray number-format style sticks to the cell.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "61.833.44"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.42%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.419.93"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.46%  "
Set-TextValue $ws.Cells.Item(4, 5) "  -0.12%  "
Set-TextValue $ws.Cells.Item(5, 4) "409.99"
Set-TextValue $ws.Cells.Item(5, 5) "  +0.20%  "
Set-TextValue $ws.Cells.Item(6, 4) "129.19"
Set-TextValue $ws.Cells.Item(6, 5) "  -0.26%  "
Set-TextValue $ws.Cells.Item(7, 5) "  -0.05%  "
Set-TextValue $ws.Cells.Item(8, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.731"
Set-TextValue $ws.Cells.Item(9, 5) "  -3.60%  "
Set-TextValue $ws.Cells.Item(10, 5) "  -0.87%  "
Set-TextValue $ws.Cells.Item(11, 4) "43.36"
Set-TextValue $ws.Cells.Item(11, 5) "  +0.05%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.0000224"
Set-TextValue $ws.Cells.Item(12, 5) "  +16.11%  "
Set-TextValue $ws.Cells.Item(13, 4) "9.27"
Set-TextValue $ws.Cells.Item(13, 5) "  +5.23%  "
Set-TextValue $ws.Cells.Item(14, 4) "3.957.70"
Set-TextValue $ws.Cells.Item(14, 5) "  -0.64%  "
Set-TextValue $ws.Cells.Item(15, 5) "  +0.40%  "
Set-TextValue $ws.Cells.Item(16, 4) "21.17"
Set-TextValue $ws.Cells.Item(16, 5) "  +4.27%  "
Set-TextValue $ws.Cells.Item(17, 4) "3.415.95"
Set-TextValue $ws.Cells.Item(17, 5) "  +0.87%  "
Set-TextValue $ws.Cells.Item(18, 4) "12.35"
Set-TextValue $ws.Cells.Item(18, 5) "  +8.46%  "
Set-TextValue $ws.Cells.Item(19, 5) "  +3.29%  "
Set-TextValue $ws.Cells.Item(20, 4) "61.819.46"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.45%  "
Set-TextValue $ws.Cells.Item(21, 4) "487.51"
Set-TextValue $ws.Cells.Item(21, 5) "  +30.12%  "
Set-TextValue $ws.Cells.Item(22, 4) "91.70"
Set-TextValue $ws.Cells.Item(22, 5) "  +4.64%  "
Set-TextValue $ws.Cells.Item(23, 4) "3.33"
Set-TextValue $ws.Cells.Item(23, 5) "  +4.67%  "
Set-TextValue $ws.Cells.Item(24, 4) "13.52"
Set-TextValue $ws.Cells.Item(24, 5) "  +1.14%  "
Set-TextValue $ws.Cells.Item(25, 4) "3.32"
Set-TextValue $ws.Cells.Item(25, 5) "  +3.44%  "
Set-TextValue $ws.Cells.Item(26, 4) "34.49"
Set-TextValue $ws.Cells.Item(26, 5) "  +9.04%  "
Set-TextValue $ws.Cells.Item(27, 4) "9.28"
Set-TextValue $ws.Cells.Item(27, 5) "  +9.61%  "
Set-TextValue $ws.Cells.Item(28, 5) "  -1.13%  "
Set-TextValue $ws.Cells.Item(29, 4) "12.12"
Set-TextValue $ws.Cells.Item(29, 5) "  +2.39%  "
Set-TextValue $ws.Cells.Item(30, 5) "  -1.88%  "
Set-TextValue $ws.Cells.Item(31, 5) "  -1.62%  "
Set-TextValue $ws.Cells.Item(32, 5) "  -2.08%  "
Set-TextValue $ws.Cells.Item(33, 4) "41.99"
Set-TextValue $ws.Cells.Item(33, 5) "  -4.49%  "
Set-TextValue $ws.Cells.Item(34, 5) "  +0.06%  "
Set-TextValue $ws.Cells.Item(35, 4) "59.07"
Set-TextValue $ws.Cells.Item(35, 5) "  +12.91%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.0498"
Set-TextValue $ws.Cells.Item(36, 5) "  +0.90%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.997"
Set-TextValue $ws.Cells.Item(38, 4) "3.47"
Set-TextValue $ws.Cells.Item(38, 5) "  +3.21%  "
Set-TextValue $ws.Cells.Item(39, 5) "  +3.37%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.75"
Set-TextValue $ws.Cells.Item(40, 5) "  +17.98%  "
Set-TextValue $ws.Cells.Item(41, 4) "146.24"
Set-TextValue $ws.Cells.Item(41, 5) "  +2.07%  "
Set-TextValue $ws.Cells.Item(42, 4) "2.93"
Set-TextValue $ws.Cells.Item(42, 5) "  +0.57%  "
Set-TextValue $ws.Cells.Item(43, 5) "  +1.62%  "
Set-TextValue $ws.Cells.Item(44, 5) "  +5.67%  "
Set-TextValue $ws.Cells.Item(45, 4) "4.37"
Set-TextValue $ws.Cells.Item(45, 5) "  +9.15%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.143"
Set-TextValue $ws.Cells.Item(50, 5) "  +17.14%  "
Set-TextValue $ws.Cells.Item(51, 4) "2.135.53"
Set-TextValue $ws.Cells.Item(51, 5) "  +1.09%  "

# Rows 46-49 were re-sorted/re-ranked: coin name, link, price and volume all change
Set-TextValue $ws.Cells.Item(46, 2) "ThetaToken"
Set-TextValue $ws.Cells.Item(46, 3) "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Cells.Item(46, 4) "2.37"
Set-TextValue $ws.Cells.Item(46, 5) "  +22.20%  "
Set-TextValue $ws.Cells.Item(47, 2) "Celestia"
Set-TextValue $ws.Cells.Item(47, 3) "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Cells.Item(47, 4) "16.73"
Set-TextValue $ws.Cells.Item(47, 5) "  +0.13%  "
Set-TextValue $ws.Cells.Item(48, 2) "EnergySwap"
Set-TextValue $ws.Cells.Item(48, 3) "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(48, 4) "23.07"
Set-TextValue $ws.Cells.Item(48, 5) "  +6.00%  "
Set-TextValue $ws.Cells.Item(49, 2) "BitcoinSV"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Cells.Item(49, 4) "117.77"
Set-TextValue $ws.Cells.Item(49, 5) "  +26.80%  "
